# "test: update test excel files to match global row rules"
#
# The "Skill" sheet only had a single leading dummy row (the
# "# Dummy Row (Ignored by ExcelBinder)" marker that ExcelBinder skips),
# while "Item" effectively reserves that same leading convention. Bring
# "Skill" in line by inserting two more dummy rows directly below the
# existing one, pushing the header/data rows down accordingly.

$wb = $excel.ActiveWorkbook

$wsSkill = $wb.Worksheets.Item("Skill")
$wsSkill.Rows("2:3").Insert()
$wsSkill.Cells.Item(2, 1).Value = "# Dummy Row (Ignored by ExcelBinder)"
$wsSkill.Cells.Item(3, 1).Value = "# Dummy Row (Ignored by ExcelBinder)"

# Leave the same selection state behind on "Item" (rows 2:3 selected).
$wsItem = $wb.Worksheets.Item("Item")
$wsItem.Activate()
[void]$wsItem.Range("A2:XFD3").Select()
